# Applies scheduled-runner market-price updates to the Pandaemonium Profits workbook.
# For each affected Leve row, updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) to the freshly pulled values; clears columns that no longer apply.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3930.5
$ws.Range("I62").Value = 2100.7144
$ws.Range("J62").Value = 8200
$ws.Range("K62").Value = 2100.7144
$ws.Range("L62").Value = 8200
$ws.Range("M62").Value = -1476.7144
$ws.Range("N62").Value = -9448

$ws.Range("H65").Value = 3930.5
$ws.Range("I65").Value = 2100.7144
$ws.Range("J65").Value = 8200
$ws.Range("K65").Value = 10503.572
$ws.Range("L65").Value = 41000
$ws.Range("M65").Value = -7383.572
$ws.Range("N65").Value = -47240

$ws.Range("H96").Value = 498
$ws.Range("I96").Value = 498
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1494
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -121

$ws.Range("H100").Value = 1598.9667
$ws.Range("I100").Value = 1338.1428
$ws.Range("J100").Value = 2207.5557
$ws.Range("K100").Value = 1338.1428
$ws.Range("L100").Value = 2207.5557
$ws.Range("M100").Value = -797.1428000000001
$ws.Range("N100").Value = -3289.5557

$ws.Range("H103").Value = 863.2105
$ws.Range("I103").Value = 683.4
$ws.Range("J103").Value = 927.4286
$ws.Range("K103").Value = 2050.2
$ws.Range("L103").Value = 2782.2858
$ws.Range("M103").Value = -1464.2
$ws.Range("N103").Value = -3954.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3134.25
$ws.Range("I110").Value = 1551.375
$ws.Range("J110").Value = 4717.125
$ws.Range("K110").Value = 1551.375
$ws.Range("L110").Value = 4717.125
$ws.Range("M110").Value = 493.625
$ws.Range("N110").Value = -8807.125

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H124").Value = 30000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 30000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 1080
$ws.Range("I10").Value = 1080
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1080
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -940

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5344

$ws.Range("H112").Value = 38234.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 38234.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 38234.5
$ws.Range("N112").Value = -41188.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1170.8572
$ws.Range("I99").Value = 1066
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1066
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = 432
$ws.Range("N99").Value = -4796

$ws.Range("H122").Value = 4692.8945
$ws.Range("I122").Value = 5767.9287
$ws.Range("J122").Value = 1682.8
$ws.Range("K122").Value = 17303.7861
$ws.Range("L122").Value = 5048.4
$ws.Range("M122").Value = -14853.7861
$ws.Range("N122").Value = -9948.4

$ws.Range("H124").Value = 32333.334
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 32333.334
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 32333.334
$ws.Range("N124").Value = -37243.334

$ws.Range("H126").Value = 1170.8572
$ws.Range("I126").Value = 1066
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 3198
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -728
$ws.Range("N126").Value = -10340

$ws.Range("H134").Value = 2269.8572
$ws.Range("I134").Value = 2055.0881
$ws.Range("J134").Value = 3182.625
$ws.Range("K134").Value = 6165.2643
$ws.Range("L134").Value = 9547.875
$ws.Range("M134").Value = -3630.2643
$ws.Range("N134").Value = -14617.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2927.2727
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2885.7144
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 8657.143199999999
$ws.Range("M80").Value = -8064
$ws.Range("N80").Value = -10529.1432

$ws.Range("H83").Value = 2927.2727
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2885.7144
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 25971.4296
$ws.Range("M83").Value = -22320
$ws.Range("N83").Value = -35331.4296

$ws.Range("H112").Value = 2712.8333
$ws.Range("I112").Value = 1200
$ws.Range("J112").Value = 2778.6086
$ws.Range("K112").Value = 3600
$ws.Range("L112").Value = 8335.825800000001
$ws.Range("M112").Value = -2492
$ws.Range("N112").Value = -10551.8258

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2807.9092
$ws.Range("I126").Value = 1992.4
$ws.Range("J126").Value = 3487.5
$ws.Range("K126").Value = 5977.200000000001
$ws.Range("L126").Value = 10462.5
$ws.Range("M126").Value = -3507.200000000001
$ws.Range("N126").Value = -15402.5

$ws.Range("H132").Value = 11703.462
$ws.Range("I132").Value = 4741.625
$ws.Range("J132").Value = 22842.4
$ws.Range("K132").Value = 14224.875
$ws.Range("L132").Value = 68527.20000000001
$ws.Range("M132").Value = -11694.875
$ws.Range("N132").Value = -73587.20000000001

$ws.Range("H141").Value = 34017
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 34017
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 34017
$ws.Range("N141").Value = -44377

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 36002
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 36002
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 36002
$ws.Range("N25").Value = -36462

$ws.Range("H46").Value = 1043.7778
$ws.Range("I46").Value = 1399
$ws.Range("J46").Value = 866.1667
$ws.Range("K46").Value = 1399
$ws.Range("L46").Value = 866.1667
$ws.Range("M46").Value = -1211
$ws.Range("N46").Value = -1242.1667

$ws.Range("H82").Value = 1250
$ws.Range("I82").Value = 833.3333
$ws.Range("J82").Value = 1406.25
$ws.Range("K82").Value = 833.3333
$ws.Range("L82").Value = 1406.25
$ws.Range("M82").Value = -472.3333
$ws.Range("N82").Value = -2128.25

$ws.Range("H85").Value = 1250
$ws.Range("I85").Value = 833.3333
$ws.Range("J85").Value = 1406.25
$ws.Range("K85").Value = 833.3333
$ws.Range("L85").Value = 1406.25
$ws.Range("M85").Value = 414.6667
$ws.Range("N85").Value = -3902.25

$ws.Range("H132").Value = 3115.3333
$ws.Range("I132").Value = 2659.8386
$ws.Range("J132").Value = 4880.375
$ws.Range("K132").Value = 7979.5158
$ws.Range("L132").Value = 14641.125
$ws.Range("M132").Value = -5449.5158
$ws.Range("N132").Value = -19701.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 3941.5
$ws.Range("I3").Value = 3941.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3941.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3827.5
$ws.Range("N3").ClearContents()

$ws.Range("H18").Value = 40003.5
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 40003.5
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 40003.5
$ws.Range("N18").Value = -40349.5
$ws.Range("M18").ClearContents()

$ws.Range("H122").Value = 1246.25
$ws.Range("I122").Value = 1246.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3738.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1288.75
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1896.841
$ws.Range("I132").Value = 1098.6571
$ws.Range("J132").Value = 5000.8887
$ws.Range("K132").Value = 3295.9713
$ws.Range("L132").Value = 15002.6661
$ws.Range("M132").Value = -765.9712999999997
$ws.Range("N132").Value = -20062.6661
